$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-08-22 Friday" "2025-08-23 Saturday"

Replace-Text "406÷5=" "569÷4="
Replace-Text "810÷3=" "965÷7="
Replace-Text "549÷8=" "449÷7="
Replace-Text "212÷6=" "434÷4="
Replace-Text "917÷9=" "268÷8="

Replace-Text "273÷4=" "285÷4="
Replace-Text "672÷2=" "999÷6="
Replace-Text "707÷8=" "435÷9="
Replace-Text "363÷8=" "981÷9="
Replace-Text "311÷7=" "907÷2="

Replace-Text "625÷4=" "520÷3="
Replace-Text "991÷7=" "105÷6="
Replace-Text "295÷8=" "990÷8="
Replace-Text "858÷4=" "172÷9="
Replace-Text "633÷5=" "955÷9="

Replace-Text "803÷5=" "791÷5="
Replace-Text "144÷4=" "551÷6="
Replace-Text "479÷8=" "853÷3="
Replace-Text "480÷2=" "702÷5="
Replace-Text "607÷9=" "552÷3="

Replace-Text "321÷4=" "259÷3="
Replace-Text "168÷5=" "697÷2="
Replace-Text "202÷8=" "195÷7="
Replace-Text "108÷5=" "861÷9="
Replace-Text "966÷4=" "128÷9="
